# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
$updates = @{ 2 = 753; 3 = 4151; 4 = 117; 5 = 755 }

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
